$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the G/H (Hop/Htr_w) values between row 2 and row 3, and renumber
# column A (ID_Building) as 1,2,3 so the ranking reflects efficiency
# (1 = best), per commit message.

$g2 = $ws.Range("G2").Value2
$h2 = $ws.Range("H2").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2

$ws.Range("G2").Value2 = $g3
$ws.Range("H2").Value2 = $h3
$ws.Range("G3").Value2 = $g2
$ws.Range("H3").Value2 = $h2

$ws.Range("A2").Value2 = 1
$ws.Range("A3").Value2 = 2
$ws.Range("A4").Value2 = 3

# Update the active selection to J14 as in the saved workbook
$ws.Range("J14").Select() | Out-Null
